$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Update existing D-column values (recomputed totals) ---
$ws.Range("D2").Value = 11768
$ws.Range("D3").Value = 11513
$ws.Range("D7").Value = 11872
$ws.Range("D8").Value = 11160
$ws.Range("D12").Value = 11980
$ws.Range("D13").Value = 10658
$ws.Range("D17").Value = 12024
$ws.Range("D18").Value = 10107
$ws.Range("D22").Value = 12081
$ws.Range("D23").Value = 9677
$ws.Range("D27").Value = 12123
$ws.Range("D28").Value = 9216
$ws.Range("D32").Value = 12163
$ws.Range("D33").Value = 8808
$ws.Range("D37").Value = 12201
$ws.Range("D38").Value = 8284
$ws.Range("D42").Value = 12236
$ws.Range("D43").Value = 7716
$ws.Range("D73").Value = 5128
$ws.Range("D74").Value = 7233
$ws.Range("D76").Value = 5816

# --- Append new weekly data block: YearWeekIso 202503, LastDayOfWeek 2025-01-19 ---
# Copy date formatting (numFmtId 14) from the last existing date cell so no new style is added.
$ws.Range("B76").Copy()
$ws.Range("B77:B81").PasteSpecial(-4122)

$ws.Range("A77").Value = 202503
$ws.Range("B77").Value = 45676
$ws.Range("C77").Value = "farms_total_count"
$ws.Range("D77").Value = 12384

$ws.Range("A78").Value = 202503
$ws.Range("B78").Value = 45676
$ws.Range("C78").Value = "farms_to_examine_count"
$ws.Range("D78").Value = 4650

$ws.Range("A79").Value = 202503
$ws.Range("B79").Value = 45676
$ws.Range("C79").Value = "farms_examined_count"
$ws.Range("D79").Value = 7734

$ws.Range("A80").Value = 202503
$ws.Range("B80").Value = 45676
$ws.Range("C80").Value = "farms_examined_positive_count"
$ws.Range("D80").Value = 1473

$ws.Range("A81").Value = 202503
$ws.Range("B81").Value = 45676
$ws.Range("C81").Value = "farms_examined_negative_count"
$ws.Range("D81").Value = 6261

# --- Update sheet view / selection to match final state ---
$ws.Range("H30:H31").Select()
